$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 92.5
$ws.Range("I6").Value = 92.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 277.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -165.5
$ws.Range("N6").ClearContents()
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H98").Value = 6920.9287
$ws.Range("I98").Value = 3427.5715
$ws.Range("J98").Value = 10414.286
$ws.Range("K98").Value = 3427.5715
$ws.Range("L98").Value = 10414.286
$ws.Range("M98").Value = -1929.5715
$ws.Range("N98").Value = -13410.286
$ws.Range("H106").Value = 3678.4375
$ws.Range("I106").Value = 1914.091
$ws.Range("K106").Value = 1914.091
$ws.Range("M106").Value = -1283.091
$ws.Range("H112").Value = 1168.3658
$ws.Range("J112").Value = 1202.641
$ws.Range("L112").Value = 3607.923
$ws.Range("N112").Value = -5823.923000000001
$ws.Range("H122").Value = 6920.9287
$ws.Range("I122").Value = 3427.5715
$ws.Range("J122").Value = 10414.286
$ws.Range("K122").Value = 10282.7145
$ws.Range("L122").Value = 31242.858
$ws.Range("M122").Value = -7832.7145
$ws.Range("N122").Value = -36142.858
$ws.Range("H129").Value = 812.2929
$ws.Range("J129").Value = 853.1319
$ws.Range("L129").Value = 2559.3957
$ws.Range("N129").Value = -12559.3957
$ws.Range("H137").Value = 1324559.9
$ws.Range("I137").Value = 2165965.2
$ws.Range("J137").Value = 2351.3572
$ws.Range("K137").Value = 6497895.600000001
$ws.Range("L137").Value = 7054.071599999999
$ws.Range("M137").Value = -6495345.600000001
$ws.Range("N137").Value = -12154.0716
$ws.Range("H138").Value = 5127.27
$ws.Range("I138").Value = 876.8421
$ws.Range("J138").Value = 6124.284
$ws.Range("K138").Value = 2630.5263
$ws.Range("L138").Value = 18372.852
$ws.Range("M138").Value = 2509.4737
$ws.Range("N138").Value = -28652.852
$ws.Range("H141").Value = 6439.6
$ws.Range("I141").Value = 7111.5454
$ws.Range("J141").Value = 3271.8572
$ws.Range("K141").Value = 21334.6362
$ws.Range("L141").Value = 9815.571599999999
$ws.Range("M141").Value = -16154.6362
$ws.Range("N141").Value = -20175.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 878
$ws.Range("I2").Value = 674.5
$ws.Range("J2").Value = 1963.3334
$ws.Range("K2").Value = 674.5
$ws.Range("L2").Value = 1963.3334
$ws.Range("M2").Value = -561.5
$ws.Range("N2").Value = -2189.3334
$ws.Range("H32").Value = 4358.6553
$ws.Range("I32").Value = 3769.3396
$ws.Range("K32").Value = 3769.3396
$ws.Range("M32").Value = -3482.3396
$ws.Range("H74").Value = 4046.879
$ws.Range("I74").Value = 4853.25
$ws.Range("J74").Value = 2806.3076
$ws.Range("K74").Value = 4853.25
$ws.Range("L74").Value = 2806.3076
$ws.Range("M74").Value = -3979.25
$ws.Range("N74").Value = -4554.3076
$ws.Range("H77").Value = 4046.879
$ws.Range("I77").Value = 4853.25
$ws.Range("J77").Value = 2806.3076
$ws.Range("K77").Value = 24266.25
$ws.Range("L77").Value = 14031.538
$ws.Range("M77").Value = -19898.25
$ws.Range("N77").Value = -22767.538
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 878
$ws.Range("I116").Value = 674.5
$ws.Range("J116").Value = 1963.3334
$ws.Range("K116").Value = 674.5
$ws.Range("L116").Value = 1963.3334
$ws.Range("M116").Value = 1619.5
$ws.Range("N116").Value = -6551.3334
$ws.Range("H132").Value = 2268.44
$ws.Range("I132").Value = 1332.2106
$ws.Range("K132").Value = 3996.6318
$ws.Range("M132").Value = -1466.6318

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 878
$ws.Range("I3").Value = 674.5
$ws.Range("J3").Value = 1963.3334
$ws.Range("K3").Value = 674.5
$ws.Range("L3").Value = 1963.3334
$ws.Range("M3").Value = -560.5
$ws.Range("N3").Value = -2191.3334
$ws.Range("H7").Value = 20122.312
$ws.Range("I7").Value = 17437.4
$ws.Range("J7").Value = 21342.727
$ws.Range("K7").Value = 17437.4
$ws.Range("L7").Value = 21342.727
$ws.Range("M7").Value = -17324.4
$ws.Range("N7").Value = -21568.727
$ws.Range("H107").Value = 1787
$ws.Range("I107").Value = 1794.2106
$ws.Range("J107").Value = 1769.875
$ws.Range("K107").Value = 1794.2106
$ws.Range("L107").Value = 1769.875
$ws.Range("M107").Value = 125.7893999999999
$ws.Range("N107").Value = -5609.875
$ws.Range("H134").Value = 2327.814
$ws.Range("I134").Value = 1783.6666
$ws.Range("J134").Value = 3246.0625
$ws.Range("K134").Value = 5350.9998
$ws.Range("L134").Value = 9738.1875
$ws.Range("M134").Value = -2815.9998
$ws.Range("N134").Value = -14808.1875
$ws.Range("H135").Value = 45116.668
$ws.Range("J135").Value = 45116.668
$ws.Range("L135").Value = 45116.668
$ws.Range("N135").Value = -55256.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7408580
$ws.Range("I16").Value = 8548053
$ws.Range("J16").Value = 2006
$ws.Range("K16").Value = 8548053
$ws.Range("L16").Value = 2006
$ws.Range("M16").Value = -8547766
$ws.Range("N16").Value = -2580
$ws.Range("H31").Value = 3280.8696
$ws.Range("I31").Value = 1215.375
$ws.Range("J31").Value = 8002
$ws.Range("K31").Value = 1215.375
$ws.Range("L31").Value = 8002
$ws.Range("M31").Value = -920.375
$ws.Range("N31").Value = -8592
$ws.Range("H34").Value = 3280.8696
$ws.Range("I34").Value = 1215.375
$ws.Range("J34").Value = 8002
$ws.Range("K34").Value = 1215.375
$ws.Range("L34").Value = 8002
$ws.Range("M34").Value = -1013.375
$ws.Range("N34").Value = -8406
$ws.Range("H113").Value = 7408580
$ws.Range("I113").Value = 8548053
$ws.Range("J113").Value = 2006
$ws.Range("K113").Value = 8548053
$ws.Range("L113").Value = 2006
$ws.Range("M113").Value = -8545883
$ws.Range("N113").Value = -6346
$ws.Range("H132").Value = 2415
$ws.Range("I132").Value = 1408.5
$ws.Range("J132").Value = 4629.3
$ws.Range("K132").Value = 4225.5
$ws.Range("L132").Value = 13887.9
$ws.Range("M132").Value = -1695.5
$ws.Range("N132").Value = -18947.9
$ws.Range("H134").Value = 5563.222
$ws.Range("I134").Value = 7800.2
$ws.Range("J134").Value = 2767
$ws.Range("K134").Value = 23400.6
$ws.Range("L134").Value = 8301
$ws.Range("M134").Value = -20865.6
$ws.Range("N134").Value = -13371

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9252
$ws.Range("J62").Value = 9252
$ws.Range("L62").Value = 27756
$ws.Range("N62").Value = -29128
$ws.Range("H65").Value = 9252
$ws.Range("J65").Value = 9252
$ws.Range("L65").Value = 83268
$ws.Range("N65").Value = -90132
$ws.Range("H113").Value = 1418.5385
$ws.Range("I113").Value = 1287.125
$ws.Range("J113").Value = 1628.8
$ws.Range("K113").Value = 3861.375
$ws.Range("L113").Value = 4886.4
$ws.Range("M113").Value = -1691.375
$ws.Range("N113").Value = -9226.4
$ws.Range("H122").Value = 2371.8518
$ws.Range("I122").Value = 868.3570999999999
$ws.Range("J122").Value = 2898.075
$ws.Range("K122").Value = 7815.2139
$ws.Range("L122").Value = 26082.675
$ws.Range("M122").Value = -5365.2139
$ws.Range("N122").Value = -30982.675
$ws.Range("H129").Value = 2343.2
$ws.Range("I129").Value = 2176.7693
$ws.Range("J129").Value = 2652.2856
$ws.Range("K129").Value = 6530.3079
$ws.Range("L129").Value = 7956.8568
$ws.Range("M129").Value = -1530.3079
$ws.Range("N129").Value = -17956.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 30000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 30000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 30000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -30226
$ws.Range("H28").Value = 30000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 30000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 30000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -30464
$ws.Range("H37").Value = 30000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 30000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 30000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -30214
$ws.Range("H40").Value = 7073.2354
$ws.Range("I40").Value = 6812.846
$ws.Range("J40").Value = 7919.5
$ws.Range("K40").Value = 6812.846
$ws.Range("L40").Value = 7919.5
$ws.Range("M40").Value = -6676.846
$ws.Range("N40").Value = -8191.5
$ws.Range("H100").Value = 3182.4546
$ws.Range("I100").Value = 1367.1666
$ws.Range("J100").Value = 5360.8
$ws.Range("K100").Value = 1367.1666
$ws.Range("L100").Value = 5360.8
$ws.Range("M100").Value = -826.1666
$ws.Range("N100").Value = -6442.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5243.773
$ws.Range("I136").Value = 2977.2
$ws.Range("J136").Value = 10100.714
$ws.Range("K136").Value = 8931.599999999999
$ws.Range("L136").Value = 30302.142
$ws.Range("M136").Value = -6381.599999999999
$ws.Range("N136").Value = -35402.142
